$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 79799.39999999999
$ws.Range("C12").Value = "'"
$ws.Range("D12").Value = 2785808.8
$ws.Range("E12").Value = 5677761.6
$ws.Range("F12").Value = 3873796.9
$ws.Range("G12").Value = 7771203.7
$ws.Range("H12").Value = 20550161.6
$ws.Range("I12").Value = 13529838.2
$ws.Range("J12").Value = 3378816.4
$ws.Range("K12").Value = 12719229.1
$ws.Range("L12").Value = "'"
$ws.Range("M12").Value = 206186292.4
$ws.Range("N12").Value = 42663432.9
$ws.Range("O12").Value = 56978366.4
$ws.Range("P12").Value = 4428452.3
$ws.Range("Q12").Value = 41634031.2
$ws.Range("R12").Value = "'"
$ws.Range("S12").Value = "'"
$ws.Range("T12").Value = 22497.3
$ws.Range("U12").Value = 292251792.7

# Reset the quote-prefix style picked up from the empty-text trick above
# so these blank cells stay on the default style (matches source row).
$ws.Range("C12").Style = "Normal"
$ws.Range("L12").Style = "Normal"
$ws.Range("R12").Style = "Normal"
$ws.Range("S12").Style = "Normal"

# A12 (the year label) should carry the same bold/bordered header style
# as the other year cells in column A.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
